$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 64.94252572409628
$ws.Range("C2").Value = 94.7083675317229
$ws.Range("D2").Value = 98.82187989694229
$ws.Range("E2").Value = 98.28123052373068
$ws.Range("F2").Value = 98.42666706849232
$ws.Range("G2").Value = 97.89243103268453
$ws.Range("H2").Value = 97.52068520966016
$ws.Range("I2").Value = 96.21924031811159

$ws.Range("B3").Value = 66.94442694611192
$ws.Range("C3").Value = 94.17750568227869
$ws.Range("D3").Value = 99.82731784166599
$ws.Range("E3").Value = 98.91690626040783
$ws.Range("F3").Value = 98.6311397405362
$ws.Range("G3").Value = 98.15307223342572
$ws.Range("H3").Value = 97.60050103818924
$ws.Range("I3").Value = 96.20732313761773

$ws.Range("B4").Value = 86.33081026460692
$ws.Range("C4").Value = 93.09258689473516
$ws.Range("D4").Value = 98.78439556107871
$ws.Range("E4").Value = 98.78839018711152
$ws.Range("F4").Value = 98.47144634178163
$ws.Range("G4").Value = 98.00599725845983
$ws.Range("H4").Value = 97.51004626892733
$ws.Range("I4").Value = 96.16055114837152

$ws.Range("B5").Value = 75.97080312273226
$ws.Range("C5").Value = 92.41733894708864
$ws.Range("D5").Value = 98.75981368270992
$ws.Range("E5").Value = 98.8890111200796
$ws.Range("F5").Value = 98.41583254386596
$ws.Range("G5").Value = 97.94311010528544
$ws.Range("H5").Value = 97.46715763101524
$ws.Range("I5").Value = 96.18455250458092

$ws.Range("B6").Value = 75.19412410012055
$ws.Range("C6").Value = 95.2242608908648
$ws.Range("D6").Value = 98.70140094385621
$ws.Range("E6").Value = 98.90394798629896
$ws.Range("F6").Value = 98.3932115870541
$ws.Range("G6").Value = 97.97034026156507
$ws.Range("H6").Value = 97.48760468379996
$ws.Range("I6").Value = 96.13076296832848
